$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $xmlPayload = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          $innerXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $r.InsertXML($xmlPayload)
}

# ---------------------------------------------------------------------------
# 1) Paragraph 1: split the run that contains "SourceTree" into three runs,
#    wrapping "SourceTree" with proofErr spell-check markers, matching the
#    target diff exactly while leaving the rest of the paragraph untouched.
# ---------------------------------------------------------------------------
$p1Xml = @"
<w:p w:rsidR="0005211A" w:rsidRDefault="00616F95" w:rsidP="00616F95">
  <w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>
  <w:r><w:t xml:space="preserve">20/07/2017 16:07 </w:t></w:r>
  <w:r><w:tab/><w:t xml:space="preserve">Files have been downloaded and added to a GitHub repository. Hosted on localhost through Node.js and managed using </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>SourceTree</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>.</w:t></w:r>
  <w:r w:rsidR="00822BB2"><w:t xml:space="preserve"> Sublime Text used for IDE.</w:t></w:r>
</w:p>
"@
Replace-ParagraphXml 1 $p1Xml

# ---------------------------------------------------------------------------
# 2) Paragraph 4 ("Replace content"): drop the _GoBack bookmark that used to
#    sit at the end of this paragraph (it moves to the new last paragraph).
# ---------------------------------------------------------------------------
$p4Xml = @"
<w:p w:rsidR="00A1623C" w:rsidRDefault="00A1623C" w:rsidP="00616F95">
  <w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>
  <w:r><w:t>20/07/2017 17:18</w:t></w:r>
  <w:r><w:tab/><w:t>Completed Basic task ‘Replace content’.</w:t></w:r>
</w:p>
"@
Replace-ParagraphXml 4 $p4Xml

# ---------------------------------------------------------------------------
# 3) Paragraph 5 (previously empty) + a brand-new paragraph 6: give the
#    empty paragraph its "Add content" text and append the new "Modify your
#    new content" paragraph (with its own proofErr-wrapped "OnClick" run and
#    the relocated _GoBack bookmark) right after it, in one InsertXML call.
# ---------------------------------------------------------------------------
$p56Xml = @"
<w:p>
  <w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>
  <w:r><w:t>20/07/2017 18:12</w:t></w:r>
  <w:r><w:tab/><w:t>Completed Basic task ‘Add content’.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr><w:ind w:left="2127" w:hanging="2127"/></w:pPr>
  <w:r><w:t>20/07/2017 19:17</w:t></w:r>
  <w:r><w:tab/><w:t xml:space="preserve">Managed to create Handlebars template for accordion element with content for each section. Came across problem where JavaScript was unable to loop through array of sections and add </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>OnClick</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> event allowing accordion to function properly.</w:t></w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
"@
Replace-ParagraphXml 5 $p56Xml

# ---------------------------------------------------------------------------
# 4) InsertXML cannot remove the document's very last paragraph mark, so the
#    step above left a trailing empty paragraph after our new content. Merge
#    it away by deleting from just before the new last paragraph's mark
#    through the end of that trailing empty paragraph.
# ---------------------------------------------------------------------------
if ($d.Paragraphs.Count -gt 6) {
    $secondToLast = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $trim = $d.Range($secondToLast.Range.End - 1, $last.Range.End)
    $trim.Delete()
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
